$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header / summary area updates
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 1597358      # VALOR MORA total
$ws.Range("C13").Value = 10           # Cant. Trabajadores
$ws.Range("F13").Value = 8            # Cant. Periodos

# ---------------------------------------------------------------------------
# 2) Make room for the new data rows.
#    Existing data occupies rows 16-32 (17 rows); the new table needs rows
#    16-39 (24 rows), i.e. 7 additional rows. Inserting those rows just
#    above the current last row (32) pushes that special "last row"
#    (bottom-border) styled row down to row 39, and shifts the blank gap +
#    footer rows (37-38) down to 44-45 - exactly where the target wants them.
# ---------------------------------------------------------------------------
$ws.Range("A32:A38").EntireRow.Insert()

# Re-apply the regular data-row formatting (copied from row 31, a normal
# data row) onto the 7 freshly inserted rows so borders/number formats
# match the rest of the table instead of Excel's generic insert defaults.
$ws.Range("B31:J31").Copy()
$ws.Range("B32:J38").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Write the new table contents (rows 16-39).
#    Column layout: B=Tipo Doc, C=N Doc Trabajador, D=Nombre Trabajador,
#    E=Periodo Mora, F=Valor Mora, G=Salario Basico
# ---------------------------------------------------------------------------
$data = @(
    @("CC","1047446431","ALFREDO JOSE ARRIETA GUERRA","2507",56940,1423500),
    @("CC","1047446431","ALFREDO JOSE ARRIETA GUERRA","2506",56940,1423500),
    @("CC","1047446431","ALFREDO JOSE ARRIETA GUERRA","2301",36341,1423500),
    @("CC","1047446431","ALFREDO JOSE ARRIETA GUERRA","2212",36341,1423500),
    @("CC","1047446431","ALFREDO JOSE ARRIETA GUERRA","2211",36341,1423500),
    @("CC","1047446431","ALFREDO JOSE ARRIETA GUERRA","2210",36341,1423500),
    @("CC","1047446431","ALFREDO JOSE ARRIETA GUERRA","2209",36341,1423500),
    @("CC","1143358596","NATALIA LEDEZMA COPETE","2507",56940,1423500),
    @("CC","1143358596","NATALIA LEDEZMA COPETE","2506",56940,1423500),
    @("CC","1044911418","EDWIN GUILLERMO PALOMINO CASTRO","2507",56940,1423500),
    @("CC","1044911418","EDWIN GUILLERMO PALOMINO CASTRO","2506",56940,1423500),
    @("CC","1054541261","FABIO HERNAN ARIAS NIETO","2507",60000,1500000),
    @("CC","1054541261","FABIO HERNAN ARIAS NIETO","2506",60000,1500000),
    @("CC","1048604236","FEDERICO QUIONES VIVANCO","2507",56940,1423500),
    @("CC","1048604236","FEDERICO QUIONES VIVANCO","2506",56940,1423500),
    @("CC","37271353","ALIX YAMILE BACCA SUAREZ","2507",140000,3500000),
    @("CC","37271353","ALIX YAMILE BACCA SUAREZ","2506",140000,3500000),
    @("CC","1090388996","EDGAR HUMBERTO BACCA SUAREZ","2507",160000,4000000),
    @("CC","1090388996","EDGAR HUMBERTO BACCA SUAREZ","2506",160000,4000000),
    @("CC","45519728","CARMEN INES RICAURTE BURGOS","2507",56940,1423500),
    @("CC","45519728","CARMEN INES RICAURTE BURGOS","2506",56940,1423500),
    @("CC","13198256","ALVARO ELIECER BACCA SUAREZ","2507",56940,1423500),
    @("CC","13198256","ALVARO ELIECER BACCA SUAREZ","2506",56940,1423500),
    @("CC","1002258091","JOSE LUIS MARMOLEJO ALVAREZ","2307",12373,1160000)
)

$r = 16
foreach ($row in $data) {
    $ws.Range("B$r").Value = $row[0]
    $ws.Range("C$r").Value = $row[1]
    $ws.Range("D$r").Value = $row[2]
    $ws.Range("E$r").Value = $row[3]
    $ws.Range("F$r").Value = $row[4]
    $ws.Range("G$r").Value = $row[5]
    $r++
}
